# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The workbook is an "Estado de Cuenta" (account statement) report for
# NIT 9013893581. This edit refreshes the summary totals and replaces the
# worker/period detail table (rows 16-32) with the new data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ---- Summary header values -------------------------------------------
# "VALOR MORA" total (E11) and the worker/period counters (C13 / F13).
$ws.Range("E11").Value = 1124682
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 17

# ---- Detail table (rows 16-32) ----------------------------------------
# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#          E=Periodo Mora, F=Valor Mora, G=Salario Basico
$rows = @(
    @{ Row = 16; Doc = "1044910398";  Nombre = "MARIA CAMILA OROZCO CUETO";    Periodo = "2106"; Valor = 36341; Salario = 908526 },
    @{ Row = 17; Doc = "1044910398";  Nombre = "MARIA CAMILA OROZCO CUETO";    Periodo = "2107"; Valor = 36341; Salario = 908526 },
    @{ Row = 18; Doc = "1044926865";  Nombre = "JESUS ALBERTO CASTRO CHAMORRO"; Periodo = "2208"; Valor = 84000; Salario = 1000000 },
    @{ Row = 19; Doc = "1044926865";  Nombre = "JESUS ALBERTO CASTRO CHAMORRO"; Periodo = "2209"; Valor = 84000; Salario = 1000000 },
    @{ Row = 20; Doc = "1044926865";  Nombre = "JESUS ALBERTO CASTRO CHAMORRO"; Periodo = "2210"; Valor = 40000; Salario = 1000000 },
    @{ Row = 21; Doc = "1002060706";  Nombre = "DANIEL BLANCO OJEDA";          Periodo = "2402"; Valor = 52000; Salario = 1300000 },
    @{ Row = 22; Doc = "1002060706";  Nombre = "DANIEL BLANCO OJEDA";          Periodo = "2403"; Valor = 52000; Salario = 1300000 },
    @{ Row = 23; Doc = "1002060706";  Nombre = "DANIEL BLANCO OJEDA";          Periodo = "2404"; Valor = 52000; Salario = 1300000 },
    @{ Row = 24; Doc = "1143378290";  Nombre = "LUIS ALFONSO PEREZ GARCIA";    Periodo = "2412"; Valor = 48000; Salario = 2000000 },
    @{ Row = 25; Doc = "1143378290";  Nombre = "LUIS ALFONSO PEREZ GARCIA";    Periodo = "2501"; Valor = 80000; Salario = 2000000 },
    @{ Row = 26; Doc = "1143378290";  Nombre = "LUIS ALFONSO PEREZ GARCIA";    Periodo = "2502"; Valor = 80000; Salario = 2000000 },
    @{ Row = 27; Doc = "1143378290";  Nombre = "LUIS ALFONSO PEREZ GARCIA";    Periodo = "2503"; Valor = 80000; Salario = 2000000 },
    @{ Row = 28; Doc = "1143378290";  Nombre = "LUIS ALFONSO PEREZ GARCIA";    Periodo = "2504"; Valor = 80000; Salario = 2000000 },
    @{ Row = 29; Doc = "1143378290";  Nombre = "LUIS ALFONSO PEREZ GARCIA";    Periodo = "2505"; Valor = 80000; Salario = 2000000 },
    @{ Row = 30; Doc = "1143378290";  Nombre = "LUIS ALFONSO PEREZ GARCIA";    Periodo = "2506"; Valor = 80000; Salario = 2000000 },
    @{ Row = 31; Doc = "1143378290";  Nombre = "LUIS ALFONSO PEREZ GARCIA";    Periodo = "2507"; Valor = 80000; Salario = 2000000 },
    @{ Row = 32; Doc = "1143378290";  Nombre = "LUIS ALFONSO PEREZ GARCIA";    Periodo = "2508"; Valor = 80000; Salario = 2000000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 3).Value = $r.Doc      # C - N Doc Trabajador
    $ws.Cells.Item($n, 4).Value = $r.Nombre   # D - Nombre Trabajador
    $ws.Cells.Item($n, 5).Value = $r.Periodo  # E - Periodo Mora
    $ws.Cells.Item($n, 6).Value = $r.Valor    # F - Valor Mora
    $ws.Cells.Item($n, 7).Value = $r.Salario  # G - Salario Basico
}
